$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.659.79"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.644.07"
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("E6").Value = "  +1.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.01"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  +0.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0627"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.68%  "
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "1.873.04"
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("D13").Value = "1.661.76"
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("E14").Value = "  +2.23%  "
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.99%  "
$ws.Range("D17").Value = "26.720.53"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.49%  "
$ws.Range("E23").Value = "  +2.17%  "
$ws.Range("E24").Value = "  +14.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("E26").Value = "  +0.34%  "
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("E28").Value = "  +5.17%  "
$ws.Range("E29").Value = "  +1.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0518"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.21%  "
$ws.Range("E31").Value = "  +0.99%  "
$ws.Range("E32").Value = "  +3.05%  "
$ws.Range("E33").Value = "  +3.18%  "
$ws.Range("D34").Value = "1.275.25"
$ws.Range("E34").Value = "  +4.74%  "
$ws.Range("E35").Value = "  +3.12%  "
$ws.Range("E36").Value = "  +5.20%  "
$ws.Range("E37").Value = "  +0.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.536"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.50%  "
$ws.Range("E39").Value = "  +3.77%  "
$ws.Range("E40").Value = "  +0.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.816"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.80%  "
$ws.Range("E43").Value = "  +2.56%  "
$ws.Range("D44").Value = "1.783.42"
$ws.Range("E44").Value = "  +1.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.82"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.71%  "
$ws.Range("E47").Value = "  +1.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0517"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.99%  "
$ws.Range("E50").Value = "  +3.22%  "
$ws.Range("E51").Value = "  -0.63%  "
